$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells are treated as text so values like "308.95", "0.1807",
# "1,660.58%" etc. round-trip exactly as strings instead of being reinterpreted
# as numbers/percentages by Excel.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '308.95'
$ws.Range("E2").Value = '0.54%'

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '41.27'
$ws.Range("E3").Value = '2.81%'

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '5.121'
$ws.Range("E4").Value = '0.71%'

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07649'
$ws.Range("E5").Value = '-0.21%'

$ws.Range("B6:E6").NumberFormat = "@"
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").Value = '1.619'
$ws.Range("E6").Value = '0.40%'

$ws.Range("B7:E7").NumberFormat = "@"
$ws.Range("B7").Value = 'BTSEToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D7").Value = '2.473'
$ws.Range("E7").Value = '2.21%'

$ws.Range("B8:E8").NumberFormat = "@"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9092'
$ws.Range("E8").Value = '-0.64%'

$ws.Range("B9:E9").NumberFormat = "@"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1184'
$ws.Range("E9").Value = '13.78%'

$ws.Range("B10:E10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1807'
$ws.Range("E10").Value = '1.50%'

$ws.Range("B11:E11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.09171'
$ws.Range("E11").Value = '-2.59%'

$ws.Range("B12:E12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.04251'
$ws.Range("E12").Value = '-4.28%'

$ws.Range("B13:E13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.1042'
$ws.Range("E13").Value = '-1.27%'

$ws.Range("B14:E14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001250'
$ws.Range("E14").Value = '-0.99%'

$ws.Range("B15:E15").NumberFormat = "@"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005874'
$ws.Range("E15").Value = '0.90%'

$ws.Range("B16:E16").NumberFormat = "@"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.356'
$ws.Range("E16").Value = '-0.11%'

$ws.Range("B17:E17").NumberFormat = "@"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '4.276'
$ws.Range("E17").Value = '0.65%'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.44%'

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '6.901'
$ws.Range("E19").Value = '-0.61%'

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1368'
$ws.Range("E20").Value = '2.06%'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.88%'

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04052'
$ws.Range("E22").Value = '-2.32%'

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001278'
$ws.Range("E23").Value = '5.99%'

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004022'
$ws.Range("E24").Value = '-2.08%'

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001271'
$ws.Range("E25").Value = '-2.34%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003744'

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02435'
$ws.Range("E38").Value = '-0.95%'

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05267'
$ws.Range("E39").Value = '1.61%'

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007806'
$ws.Range("E40").Value = '-1.32%'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.03%'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-4.33%'

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001950'
$ws.Range("E43").Value = '-0.05%'

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007553'
$ws.Range("E44").Value = '1.80%'

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3077'
$ws.Range("E45").Value = '0.36%'

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006885'
$ws.Range("E46").Value = '6.81%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.12%'

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '0.08006'
$ws.Range("E48").Value = '1,660.58%'

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003000'
$ws.Range("E49").Value = '-0.13%'

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '-0.12%'

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '-0.12%'
